# Adding more code for plotting
# - Adds a new "Sheet2" worksheet (placed after "Sheet1") containing
#   Length / Yield sample data used for a plot.
# - Leaves Sheet1's data untouched, but updates its selection and makes
#   Sheet2 the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet and put it right after Sheet1 -----------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet2"

# Re-fetch live references (the handles captured before the structural
# edits above can go stale) and reorder: Sheet1, then Sheet2.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2.Move($null, $ws1)

# Re-fetch again now that the sheet collection has settled.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Populate Sheet2 ----------------------------------------------------
$ws2.Range("B1").Value = "Length"
$ws2.Range("C1").Value = "Yield"

$length = @(1.45, 2.32, 3.12, 4.13, 5.67, 6.86, 7.52, 8.32, 9.12, 10.34)
$yield  = @(2.1, 4.12, 6.88, 8.12, 10.43, 11.45, 15.07, 16.56, 17.86, 21.09)

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws2.Range("A$r").Value = $i
    $ws2.Range("B$r").Value = $length[$i]
    $ws2.Range("C$r").Value = $yield[$i]
}

# --- View/selection state ------------------------------------------------
# Sheet1 keeps its data, but the selection moves to A1:D4 and it is no
# longer the visible/active tab.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("A1:D4").Select()

# Sheet2 becomes the active sheet, selected at C2.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("C2").Select()
$ws2.Activate()
